# jn: revisión de tipos de variables, ejemplos y definiciones
#
# Updates the "Clasificación" column of the 02_1_diccionario sheet with the
# reviewed classification labels, and refreshes the active-cell selections
# on Hoja2 and 02_1_diccionario left over from that review pass.

$wb = $excel.ActiveWorkbook

# --- 02_1_diccionario: revised "Clasificación" values (column I) ---------
$dic = $wb.Worksheets.Item("02_1_diccionario")

$dic.Range("I2").Value = "Identificación"
$dic.Range("I3").Value = "Temporalidad"
$dic.Range("I4").Value = "Localización"
$dic.Range("I5").Value = "Localización"
$dic.Range("I6").Value = "Descripción"
$dic.Range("I7").Value = "Transacción"
$dic.Range("I8").Value = "Descripción"

# --- Hoja2: selection left on this sheet after the review -----------------
$hoja2 = $wb.Worksheets.Item("Hoja2")
[void]$hoja2.Activate()
$hoja2.Range("C6").Select() | Out-Null

# 02_1_diccionario stays the active/front sheet, with its own scroll and
# selection updated (re-activate it last so it keeps the focus).
[void]$dic.Activate()
$dic.Range("C7").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 4
